$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[49.813737449918044, 50.06823071872015]"
$ws.Range("U2").Value = "[49.99485130692839, 50.169882464036185]"

$ws.Range("M3").Value = "[49.93364380700902, 50.21019661503114]"
$ws.Range("U3").Value = "[49.971461934704706, 50.12333886070036]"

$ws.Range("M4").Value = "[49.865411375166424, 50.19096391862869]"
$ws.Range("U4").Value = "[49.889935361648625, 50.066595672379606]"

$ws.Range("M5").Value = "[49.82575597582388, 50.14814143514569]"
$ws.Range("U5").Value = "[49.9885351504619, 50.16556565269782]"

$ws.Range("M6").Value = "[49.92486771126618, 50.16539913996673]"
$ws.Range("U6").Value = "[49.94899564660623, 50.099490212167204]"

$ws.Range("M7").Value = "[49.93640770118993, 50.173240851649986]"
$ws.Range("U7").Value = "[49.989318430417, 50.152699680545155]"
